$d = $word.ActiveDocument

# The table definition for "PrestaServiço" lists its columns as:
#   (id_prestador, id_imovel, data_inicio, data_fim, valor_mensal)
# and the commit renames the "valor_mensal" column to just "valor".
# In the source OOXML the phrase ", data_inicio, data_fim, " and the
# word "valor_mensal" are two separate (but identically-formatted) runs
# that must be collapsed into a single run reading
#   ", data_inicio, data_fim, valor"
# while the run(s) that follow (the closing ")" and the space after it)
# must stay exactly as they were.

$oldPhrase = ", data_inicio, data_fim, valor_mensal"
$newPhrase = ", data_inicio, data_fim, valor"

# Locate the phrase spanning the two runs that need to be merged.
$target = $d.Content
$find = $target.Find
$find.ClearFormatting()
$find.Text = $oldPhrase
$found = $target.Find.Execute()
if (-not $found) {
    throw "Could not find the text '$oldPhrase' to update"
}

# The COM runtime coalesces a run with every following run that shares
# its formatting whenever the run's text is edited. That would also
# swallow the ")" (and trailing space) that immediately follow
# "valor_mensal" even though the diff leaves them untouched. To stop
# the coalescing at the right place we briefly give the very next
# character a different (Bold) formatting so it no longer matches the
# run being edited, perform the text replacement, then restore its
# original formatting.
$guard = $d.Range($target.End, $target.End + 1)
$guard.Bold = $true

$target.Text = $newPhrase

$guardNewStart = $target.Start + $newPhrase.Length
$guardRestored = $d.Range($guardNewStart, $guardNewStart + 1)
$guardRestored.Bold = $false
